# _BD (Autosaved).xlsx — "Ready con Base de Datos" update
# Applies the grade-input corrections (columns H/I/J/M) for the affected
# students, lets the dependent formulas (D/E/F/G/N) recalculate, extends
# the AutoFilter / _FilterDatabase range down through row 99, and updates
# the saved view state (frozen-pane scroll position + active selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Grade-input corrections per row (only columns H, I, J, K, L, M are
#    raw inputs; D, E, F, G and N are formulas and recalc automatically).
# ---------------------------------------------------------------------

# Rows where the "Exp.20" input (column I) was cleared out entirely.
$clearedI = @(3, 4, 8, 18, 23, 24, 25, 26, 49, 95, 99)
foreach ($r in $clearedI) {
    $ws.Range("I$r").Value = $null
}

# Row 9 also had its "AS" input (column H) explicitly set to 0.
$ws.Range("H9").Value = 0

# New/updated "PF.30" input (column M) values.
$mValues = @{
    6  = 28
    9  = 18
    12 = 29
    14 = 7
    16 = 30
    19 = 29
    20 = 29
    27 = 1
    28 = 15
    29 = 23
    30 = 27
    31 = 30
    32 = 1
    35 = 28
    42 = 30
    45 = 30
    47 = 30
    50 = 30
    52 = 9
    65 = 30
    71 = 30
    72 = 30
    74 = 23
    78 = 28
    79 = 28
    80 = 30
    81 = 30
    83 = 30
    85 = 30
    87 = 30
    90 = 28
    94 = 23
}
foreach ($r in $mValues.Keys) {
    $ws.Range("M$r").Value = $mValues[$r]
}

# ---------------------------------------------------------------------
# 2. Extend the AutoFilter range from A1:O1 to A1:O99 (and keep the
#    _xlnm._FilterDatabase defined name in sync with it).
# ---------------------------------------------------------------------

$ws.AutoFilterMode = $false
$ws.Range("A1:O99").AutoFilter()

foreach ($i in 1..$wb.Names.Count()) {
    $n = $wb.Names.Item($i)
    if ($n.Name() -like "*_FilterDatabase") {
        $n.RefersTo = "='Page 1'!`$A`$1:`$O`$99"
    }
}

# ---------------------------------------------------------------------
# 3. Update the saved view state: select I99 (last edited cell) so the
#    sheet's stored selection/scroll position matches the authored file.
# ---------------------------------------------------------------------

$ws.Activate()
$ws.Range("I99").Select()
